$wb = $excel.ActiveWorkbook

# Update values on "展览" sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 156
$wsExpo.Range("F4").Value = 263
$wsExpo.Range("F5").Value = 4011
$wsExpo.Range("F6").Value = 35
$wsExpo.Range("F7").Value = 445

# Update values on "全部类型" sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 156
$wsAll.Range("F4").Value = 263
$wsAll.Range("F5").Value = 4011
$wsAll.Range("F8").Value = 35
$wsAll.Range("F9").Value = 445
